$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 1500
$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 5250

$ws.Range("E3").ClearContents()

$ws.Range("E2").Select()
